# Update iteration 0 and iteration 1 rows with new benchmark values,
# and remove the iteration 2 row entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New values for row 2 (iteration 0)
$ws.Range("B2").Value = 1437
$ws.Range("C2").Value = 395
$ws.Range("D2").Value = 177
$ws.Range("E2").Value = 27
$ws.Range("F2").Value = 2722
$ws.Range("G2").Value = 4351
$ws.Range("H2").Value = 5001

# New values for row 3 (iteration 1)
$ws.Range("B3").Value = 1444
$ws.Range("C3").Value = 396
$ws.Range("D3").Value = 173
$ws.Range("E3").Value = 29
$ws.Range("F3").Value = 2654
$ws.Range("G3").Value = 4287
$ws.Range("H3").Value = 4928

# Remove the iteration 2 row (row 4) entirely, shifting cells up
$ws.Rows.Item(4).Delete()
